# "new calculations with new timeseries"
$wb = $excel.ActiveWorkbook

$cost = $wb.Worksheets.Item("cost")
$ts   = $wb.Worksheets.Item("timeseries")

# ---------------------------------------------------------------
# "cost" sheet - updated technology cost inputs
# ---------------------------------------------------------------
$cost.Range("G5").Value = 95
$cost.Range("H5").Value = 58
$cost.Range("H6").Value = 73
$cost.Range("G7").Value = 80
$cost.Range("G8").Value = 46

# ---------------------------------------------------------------
# "timeseries" sheet - updated hourly generation-mix shares
# ---------------------------------------------------------------

# Column F ("coal" share)
$ts.Range("F4").Value  = 0.67
$ts.Range("F5").Value  = 0.51
$ts.Range("F6").Value  = 0.4
$ts.Range("F7").Value  = 0.33
$ts.Range("F12").Value = 0.2
$ts.Range("F13").Value = 0.15
$ts.Range("F14").Value = 0.11
$ts.Range("F15").Value = 0.04
$ts.Range("F16").Value = 0.02
$ts.Range("F17").Value = 0.08
$ts.Range("F20").Value = 0.18
$ts.Range("F21").Value = 0.29
$ts.Range("F22").Value = 0.36
$ts.Range("F23").Value = 0.45
$ts.Range("F24").Value = 0.38
$ts.Range("F25").Value = 0.3
$ts.Range("F26").Value = 0.26

# Column G ("gas" share)
$ts.Range("G18").Value = 0.03
$ts.Range("G19").Value = 0.11
$ts.Range("G20").Value = 0.24
$ts.Range("G21").Value = 0.32
$ts.Range("G22").Value = 0.28
$ts.Range("G23").Value = 0.13
$ts.Range("G24").Value = 0.04

# Column H ("wind" share)
$ts.Range("H4").Value  = 0.51
$ts.Range("H5").Value  = 0.42
$ts.Range("H6").Value  = 0.33
$ts.Range("H7").Value  = 0.18
$ts.Range("H8").Value  = 0.15
$ts.Range("H9").Value  = 0.12
$ts.Range("H10").Value = 0.08
$ts.Range("H11").Value = 0.15
$ts.Range("H14").Value = 0.13
$ts.Range("H15").Value = 0.06
$ts.Range("H16").Value = 0.07
$ts.Range("H17").Value = 0.12
$ts.Range("H18").Value = 0.13
$ts.Range("H19").Value = 0.16
$ts.Range("H20").Value = 0.17
$ts.Range("H21").Value = 0.21
$ts.Range("H22").Value = 0.31
$ts.Range("H23").Value = 0.45
$ts.Range("H24").Value = 0.32
$ts.Range("H26").Value = 0.13
$ts.Range("H27").Value = 0.11

# Column I ("solar" share)
$ts.Range("I21").Value = 0.45
$ts.Range("I22").Value = 0.42

# Standalone inputs next to row 46
$ts.Range("N46").Value = 2200
$ts.Range("P46").Value = 1700

# ---------------------------------------------------------------
# View state: selection per sheet + which sheet/tab is active
# ---------------------------------------------------------------
$cost.Activate()
$cost.Range("H7").Select() | Out-Null

$ts.Activate()
$ts.Range("H14").Select() | Out-Null
